# Generate Report for Handoff
# Update timestamps / priority for the six "Ready for handoff" rows
# (rows 8, 9, 10, 12, 13, 14) across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-08-13 00:25:09"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-13 00:24:55"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-13 00:25:09"
}
